# Applies the metrics_6_3.xlsx edit:
#  - Row labels in column A are rotated for rows 2, 14, 26:
#       A2  : model_6_3_12 -> model_6_3_0
#       A14 : model_6_3_24 -> model_6_3_12
#       A26 : model_6_3_0  -> model_6_3_24
#  - Every data row (2..26), columns B..Q, is overwritten with the same
#    refreshed metrics vector (the ensemble's recomputed dispersion stats).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A label rotation -------------------------------------------------
$ws.Range("A2").Value  = "model_6_3_0"
$ws.Range("A14").Value = "model_6_3_12"
$ws.Range("A26").Value = "model_6_3_24"

# --- New shared metrics values for columns B..Q, applied to rows 2..26 ------
$newValues = @(
    0.5692387035740197,
    0.2355707868591003,
    -4.985891163924485,
    -0.1117018534276049,
    -0.7457095594852317,
    0.2557185134775685,
    0.4537982025430372,
    0.3460263304716791,
    0.4207359091111362,
    0.3833811197914076,
    0.2823828768330753,
    0.5056861808251918,
    0.06015717143422494,
    0.5272142884576241,
    28.72735599237364,
    44.57274171566024
)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        $col = $i + 2   # column B is index 2
        $ws.Cells.Item($row, $col).Value = $newValues[$i]
    }
}
